# update template laba rugi
# - flag both rows' "Report Bold1" (col N) / "Report Bold2" (col S) as 1 (was
#   missing / stored as boolean for one row each); now both rows carry the
#   same numeric 1 in both columns
# - widen "Account Name2" (col I / 9) and "Report Bold1" (col N / 14) columns
# - move the active selection to L6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (data row 1): add the missing Report Bold1 flag, and make Report
# Bold2 numeric (it was already 1 but stored as a boolean)
$ws.Range("N2").Value = 1
$ws.Range("S2").Value = 1

# Row 3 (data row 2): make Report Bold1 numeric (was a boolean), and add the
# missing Report Bold2 flag
$ws.Range("N3").Value = 1
$ws.Range("S3").Value = 1

# Custom column widths (Account Name2 ~13.29 chars, Report Bold1 ~21.71 chars
# in Excel's own pixel-based width model)
$ws.Columns.Item(9).ColumnWidth = 12.5
$ws.Columns.Item(14).ColumnWidth = 20.833333333333332

# Move selection
$ws.Range("L6").Select() | Out-Null
